$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Save" - reuse the same formatting as the other
# header cells (bold, centered, thin border) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data column H2:H4 ("Save" flag)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
